$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "add client" data provider: assertion columns ---
# Header for the expected-result column (S) is entered first, then the
# expected values for the two test rows, then the test data for row 3 is
# updated to use a unique client name ("amol1235" instead of "amol" so the
# "add new client" test no longer collides with an existing record), then
# the expected result for row 3, then the actual-result xpath values for
# rows 2 and 3, and finally the header for the xpath column (T).
$ws.Range("S1").Value = "expected result"
$ws.Range("S2").Value = "Client already exists!"
$ws.Range("A3").Value = "amol1235"
$ws.Range("S3").Value = "Record successfully created"
$ws.Range("T2").Value = "//div[contains(@class,'alert')]"
$ws.Range("T3").Value = "//div[contains(@class,'alert')]"
$ws.Range("T1").Value = "xpath of actual result"

# Size the new columns similarly to how Excel auto-fits them to their content.
$ws.Columns.Item(12).ColumnWidth = 11.6
$ws.Columns.Item(19).ColumnWidth = 25
$ws.Columns.Item(20).ColumnWidth = 27

# Reflect the new selection/scroll position used while authoring the data.
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("S8").Select()
